# Weekly update: two new price records were added for "Camote" (Zapallo /
# Terminal La Palmera de La Serena) ahead of the existing row 888, pushing
# the remaining 90 data rows (old 888-977) down by two rows (to 890-979).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 888 - this shifts the
# existing rows 888:977 down to 890:979 and keeps their data intact.
$ws.Rows("888:889").Insert()

# New row 888
$ws.Range("A888").Value = 8
$ws.Range("B888").Value = "Terminal La Palmera de La Serena"
$ws.Range("C888").Value = "Coquimbo"
$ws.Range("D888").Value = 45166
$ws.Range("E888").Value = 4
$ws.Range("F888").Value = 100112045
$ws.Range("G888").Value = "Zapallo"
$ws.Range("H888").Value = "Camote"
$ws.Range("I888").Value = "1a (guarda)"
$ws.Range("J888").Value = 1200
$ws.Range("K888").Value = 1100
$ws.Range("L888").Value = 1200
$ws.Range("M888").Value = 1150
$ws.Range("N888").Value = "$/kilo (volumen en unidades)"
$ws.Range("O888").Value = "Región de O'Higgins"
$ws.Range("P888").Value = 1150
$ws.Range("Q888").Value = 1
$ws.Range("R888").Value = "Hortaliza"

# New row 889
$ws.Range("A889").Value = 8
$ws.Range("B889").Value = "Terminal La Palmera de La Serena"
$ws.Range("C889").Value = "Coquimbo"
$ws.Range("D889").Value = 45166
$ws.Range("E889").Value = 4
$ws.Range("F889").Value = 100112045
$ws.Range("G889").Value = "Zapallo"
$ws.Range("H889").Value = "Camote"
$ws.Range("I889").Value = "1a nueva(o)"
$ws.Range("J889").Value = 1600
$ws.Range("K889").Value = 950
$ws.Range("L889").Value = 1000
$ws.Range("M889").Value = 975
$ws.Range("N889").Value = "$/kilo (volumen en unidades)"
$ws.Range("O889").Value = "Perú"
$ws.Range("P889").Value = 975
$ws.Range("Q889").Value = 1
$ws.Range("R889").Value = "Hortaliza"
